# Applies the "Fixed grammar and spelling errors in PMP and COO" edits to
# Concept of Operations.docx.
#
# Strategy: use Find/Replace (wdReplaceAll) against $d.Content for every
# textual fix, then relocate the stray "_GoBack" bookmark from its old
# location (an empty heading paragraph near the top of the doc) to the end
# of the very last edited sentence, matching how Word re-stamps _GoBack at
# the site of the most recent edit.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "MISSING: $find"
    }
}

# --- 1. "windows" -> "Windows" (both mentions), drop comma after "GUI" ---
Replace-Text "coding a windows desktop application" "coding a Windows desktop application"
Replace-Text "needs the windows operating system" "needs the Windows operating system"
Replace-Text "that is relevant to the user (based on user input in the GUI, and based on which user" "that is relevant to the user (based on user input in the GUI and based on which user"

# --- 2. Face-detection bullet: split sentence, drop the tracking-box clause ---
Replace-Text "it will prompt them to complete a profile for it to be able to recognize them, otherwise, it will greet them by their name and put a tracking box on their face that has a label with their name. " "it will prompt them to complete a profile for it to be able to recognize them. Otherwise, it will greet them by their name. "

# --- 3. JARVIS command bullet: split sentence, drop parens, lowercase "answer" ---
Replace-Text "it will be able to follow commands, for example (if the user starts a sentence with a keyword, it will take in a command that has already been established). The user says" "it will be able to follow commands. For example, if the user starts a sentence with a keyword, it will take in a command that has already been established. The user says"
Replace-Text "the application will then Answer " "the application will then answer "

# --- 4. Too many faces bullet ---
Replace-Text "Too many faces in the frame may cause the program to crash or get people’s profiles confused" "Too many faces in the frame may cause the program to crash or confuse people’s profiles."

# --- 5. Not understanding commands bullet ---
Replace-Text "doing something completely different to what the user requested." "doing something completely different from what the user requested."

# --- 6. Drop trailing periods on several "Needed Features" bullets ---
Replace-Text "most used applications, most visited websites." "most used applications, most visited websites"
Replace-Text "It must be able to perform speech recognition and construction." "It must be able to perform speech recognition and construction"
Replace-Text "The ability to report information from websites (such as weather)." "The ability to report information from websites (such as weather)"
Replace-Text "The ability to open and close applications, and log in or out." "The ability to open and close applications, and log in or out"

# --- 7. Mouth movement bullet ---
Replace-Text "The ability to analyze mouth movement from users to be able who its interacting with." "The ability to analyze mouth movement from users to be able to determine whom it is interacting with"

# --- 8. Animated face bullet (also strips the stray leading space) ---
Replace-Text " Have animated face talk to user." "An animated face to talk to the user"

# --- 9. Emotions/gestures bullet (strips the stray leading space) ---
Replace-Text " Recognizing emotions and gestures" "Recognition of emotions and gestures"

# --- 10. Collapse the doubled space before "users searching" ---
Replace-Text "instead of  users searching" "instead of users searching"

# --- 11. "web APIs formats" -> "web API formats" ---
Replace-Text "the difference in web APIs formats" "the difference in web API formats"

# --- 12. "couldn't" -> "cannot" ---
Replace-Text "All algorithms that couldn't be enhanced will be discussed for further research." "All algorithms that cannot be enhanced will be discussed for further research."

# --- 13. Final "Mainly, the Web APIs..." paragraph rewrite ---
Replace-Text "for the application to interpret and read for to the user, so we will have to choose whether or not to implement specific websites, although a system can be designed for add-ons and plug-ins that could increase the amount of websites supported." "for the application to interpret and present to the user, so we will have to choose whether or not to implement specific websites. However, a system can be designed for add-ons and plug-ins that could increase the amount of websites supported."

# --- 14. Relocate the "_GoBack" bookmark to the end of that last paragraph ---
try {
    $old = $d.Bookmarks("_GoBack")
    $old.Delete()
} catch {
}

$endRng = $d.Content
$found = $endRng.Find.Execute("websites supported.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $endRng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $endRng)
} else {
    Write-Output "MISSING: websites supported. (for bookmark relocation)"
}

Write-Output "done"
